# Merge the "Degeneracy of the modularity function..." FirstParagraph and the
# following BodyText quote paragraph into a single, rewritten FirstParagraph
# that explains the modularity-degeneracy issue and the Good/de
# Montjoye/Clauset (2010) citation, including inline equations.
$d = $word.ActiveDocument

$startRange = $d.Content.Duplicate
$found1 = $startRange.Find.Execute("Degeneracy of the modularity function", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "edit.ps1: could not find start anchor paragraph" }

$endRange = $d.Content.Duplicate
$found2 = $endRange.Find.Execute("Good, de Montjoye, and Clauset (2010).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "edit.ps1: could not find end anchor paragraph" }

$whole = $d.Range($startRange.Start, $endRange.End)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The first issue with modularity maximization is the the modularity function</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>Q</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">has many local optima, with similar values of</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:r><m:t>Q</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve">, but which correspond to</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">qualitatively very different partitions</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:acc><m:accPr><m:chr m:val="⃗"/></m:accPr><m:e><m:r><m:t>b</m:t></m:r></m:e></m:acc></m:oMath><w:r><w:t xml:space="preserve">. This was first reported in</w:t></w:r><w:r><w:br/></w:r><w:r><w:t xml:space="preserve">Good, de Montjoye, and Clauset (2010)</w:t></w:r><w:r><w:t xml:space="preserve">, who also show that</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><m:oMath><m:sSub><m:e><m:r><m:t>Q</m:t></m:r></m:e><m:sub><m:r><m:t>m</m:t></m:r></m:sub></m:sSub><m:r><m:t>a</m:t></m:r><m:r><m:t>x</m:t></m:r></m:oMath><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">is highly dependent on the number of</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">clusters and the size of the network, and conclude—</w:t></w:r><w:r><w:t xml:space="preserve">“</w:t></w:r><w:r><w:t xml:space="preserve">[the] modules identified</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">through modularity maximization should be treated with caution in all but the</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">most straightforward cases</w:t></w:r><w:r><w:t xml:space="preserve">”</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">(Good, de Montjoye, and Clauset 2010)</w:t></w:r><w:r><w:t xml:space="preserve">.</w:t></w:r></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$whole.InsertXML($xml)
